$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "About": update notes text to reflect the Agora (EU) data source
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Row 6 text changes in place (new wording, mentions "in some regions"/"in the US")
$wsAbout.Range("A6").Value = "Certain plant types in some regions, such as coal and natural gas in the US, are capable of running for most"

# Insert 3 new rows before the old row 12 (2 new note lines + 1 blank spacer),
# which pushes everything from the old row 12 onward down by 3 rows.
$wsAbout.Rows("12:14").Insert()

$wsAbout.Range("A12").Value = "In the EU 28 we assume a high degree of the power plants' flexibility."
$wsAbout.Range("A13").Value = "Some plant types are thus set to 1 so that they have the flexibility to bid at higher capacity factors."
# row 14 is left blank (spacer), matching the blank spacer rows elsewhere on this sheet

# ---------------------------------------------------------------------
# Sheet "BDSBaPCF": the duplicate-fuel rows (lignite, offshore wind, crude
# oil, heavy/residual fuel oil, municipal solid waste) no longer mirror
# another row via formula -- they now hold plain static values.
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("BDSBaPCF")

$wsData.Range("B13").Value = 1   # lignite            (was =B2)
$wsData.Range("B14").Value = 0   # offshore wind       (was =B6)
$wsData.Range("B15").Value = 0   # crude oil           (was =B11)
$wsData.Range("B16").Value = 0   # heavy or residual fuel oil (was =B11)
$wsData.Range("B17").Value = 1   # municipal solid waste (was =B9)
